$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 111 (shifts old rows 111-125 down to 113-127)
$ws.Range("A111:A112").EntireRow.Insert()

# New row 111: weekly "Primera" entry dated 45015
$ws.Range("A111").Value = 2
$ws.Range("B111").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C111").Value = "Coquimbo"
$ws.Range("D111").Value = 45015
$ws.Range("E111").Value = 4
$ws.Range("F111").Value = "Fruta"
$ws.Range("G111").Value = 100107
$ws.Range("H111").Value = "Otros"
$ws.Range("I111").Value = 100107011
$ws.Range("J111").Value = "Tuna"
$ws.Range("K111").Value = "Sin especificar"
$ws.Range("L111").Value = "Primera"
$ws.Range("M111").Value = 360
$ws.Range("N111").Value = 10000
$ws.Range("O111").Value = 11000
$ws.Range("P111").Value = 10500
$ws.Range("Q111").Value = "$/caja 18 kilos"
$ws.Range("R111").Value = "Provincia de Limarí"
$ws.Range("S111").Value = 583
$ws.Range("T111").Value = 18

# New row 112: weekly "Segunda" entry dated 45015
$ws.Range("A112").Value = 2
$ws.Range("B112").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C112").Value = "Coquimbo"
$ws.Range("D112").Value = 45015
$ws.Range("E112").Value = 4
$ws.Range("F112").Value = "Fruta"
$ws.Range("G112").Value = 100107
$ws.Range("H112").Value = "Otros"
$ws.Range("I112").Value = 100107011
$ws.Range("J112").Value = "Tuna"
$ws.Range("K112").Value = "Sin especificar"
$ws.Range("L112").Value = "Segunda"
$ws.Range("M112").Value = 300
$ws.Range("N112").Value = 8000
$ws.Range("O112").Value = 9000
$ws.Range("P112").Value = 8500
$ws.Range("Q112").Value = "$/caja 18 kilos"
$ws.Range("R112").Value = "Provincia de Limarí"
$ws.Range("S112").Value = 472
$ws.Range("T112").Value = 18
